# "all gems discovered up to x7" - fill in the remaining newly-discovered
# gem markers on the "Retry" sheet (the active sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully confirmed discoveries ("Q") for gems up to x7 (column K):
$ws.Range("L2").Value = "Q"   # Red / x8
$ws.Range("J4").Value = "Q"   # Yellow / x6
$ws.Range("K8").Value = "Q"   # White / x7

# Partially confirmed discovery marked with "?" (Purple / x8):
$ws.Range("L6").Value = "?"

# Move the active selection to reflect where the user left off.
$ws.Range("L5").Select() | Out-Null
